$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H changes (row 20 becomes a normal "bedrooms" row, row 26 becomes the catch row)
$ws.Range("H20").Value = "bedrooms"
$ws.Range("H26").ClearContents()

# Row 2
$ws.Range("I2").Value = "target"
$ws.Range("J2").Value = "old"
$ws.Range("K2").Value = "j"
$ws.Range("L2").Value = "stimuli/img_2pk6v.png"
$ws.Range("M2").Value = 85.08108108108108
$ws.Range("N2").Value = 66.16216216216216
$ws.Range("O2").Value = 75.62162162162161
$ws.Range("P2").Value = 37
$ws.Range("Q2").Value = 9
$ws.Range("R2").Value = 9
$ws.Range("S2").Value = 9
$ws.Range("T2").Value = 9
$ws.Range("U2").Value = 9
$ws.Range("V2").Value = 8

# Row 3
$ws.Range("I3").ClearContents()
$ws.Range("J3").Value = "new"
$ws.Range("K3").Value = "f"
$ws.Range("L3").Value = "stimuli/img_bklr1.png"
$ws.Range("M3").Value = 86.54761904761905
$ws.Range("N3").Value = 67.73809523809524
$ws.Range("O3").Value = 77.14285714285714
$ws.Range("P3").Value = 42
$ws.Range("Q3").Value = 9
$ws.Range("R3").Value = 9
$ws.Range("S3").Value = 9
$ws.Range("T3").Value = 9
$ws.Range("U3").Value = 9
$ws.Range("V3").Value = 9

# Row 4
$ws.Range("I4").Value = "target"
$ws.Range("J4").Value = "old"
$ws.Range("K4").Value = "j"
$ws.Range("L4").Value = "stimuli/img_bj2gr.png"
$ws.Range("M4").Value = 65.25
$ws.Range("N4").Value = 44.8
$ws.Range("O4").Value = 55.025
$ws.Range("P4").Value = 40
$ws.Range("Q4").Value = 4
$ws.Range("R4").Value = 4
$ws.Range("S4").Value = 4
$ws.Range("T4").Value = 4
$ws.Range("U4").Value = 4
$ws.Range("V4").Value = 4

# Row 5
$ws.Range("I5").Value = "target"
$ws.Range("J5").Value = "old"
$ws.Range("K5").Value = "j"
$ws.Range("L5").Value = "stimuli/img_jivhq.png"
$ws.Range("M5").Value = 37
$ws.Range("N5").Value = 22.26530612244898
$ws.Range("O5").Value = 29.63265306122449
$ws.Range("P5").Value = 49
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 2
$ws.Range("S5").Value = 2
$ws.Range("T5").Value = 2
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 2

# Row 6
$ws.Range("I6").Value = "target"
$ws.Range("J6").Value = "old"
$ws.Range("K6").Value = "j"
$ws.Range("L6").Value = "stimuli/img_z3yzz.png"
$ws.Range("M6").Value = 71.71052631578948
$ws.Range("N6").Value = 49.81578947368421
$ws.Range("O6").Value = 60.76315789473685
$ws.Range("P6").Value = 38
$ws.Range("Q6").Value = 5
$ws.Range("R6").Value = 5
$ws.Range("S6").Value = 5
$ws.Range("T6").Value = 5
$ws.Range("U6").Value = 5
$ws.Range("V6").Value = 5

# Row 7
$ws.Range("I7").ClearContents()
$ws.Range("J7").Value = "new"
$ws.Range("K7").Value = "f"
$ws.Range("L7").Value = "stimuli/img_6ddrx.png"
$ws.Range("M7").Value = 82.2
$ws.Range("N7").Value = 63.68571428571428
$ws.Range("O7").Value = 72.94285714285715
$ws.Range("P7").Value = 35
$ws.Range("Q7").Value = 8
$ws.Range("R7").Value = 8
$ws.Range("S7").Value = 8
$ws.Range("T7").Value = 8
$ws.Range("U7").Value = 8
$ws.Range("V7").Value = 8

# Row 8
$ws.Range("I8").ClearContents()
$ws.Range("J8").Value = "new"
$ws.Range("K8").Value = "f"
$ws.Range("L8").Value = "stimuli/img_ri0yx.png"
$ws.Range("M8").Value = 88.96969696969697
$ws.Range("N8").Value = 77.15151515151516
$ws.Range("O8").Value = 83.06060606060606
$ws.Range("P8").Value = 33
$ws.Range("Q8").Value = 10
$ws.Range("R8").Value = 10
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = 10
$ws.Range("U8").Value = 10
$ws.Range("V8").Value = 10

# Row 9
$ws.Range("I9").Value = "target"
$ws.Range("J9").Value = "old"
$ws.Range("K9").Value = "j"
$ws.Range("L9").Value = "stimuli/img_okvvw.png"
$ws.Range("M9").Value = 50.58333333333334
$ws.Range("N9").Value = 32.11111111111111
$ws.Range("O9").Value = 41.34722222222223
$ws.Range("P9").Value = 36
$ws.Range("Q9").Value = 2
$ws.Range("R9").Value = 2
$ws.Range("S9").Value = 2
$ws.Range("T9").Value = 2
$ws.Range("U9").Value = 2
$ws.Range("V9").Value = 3

# Row 10
$ws.Range("I10").Value = "target"
$ws.Range("J10").Value = "old"
$ws.Range("K10").Value = "j"
$ws.Range("L10").Value = "stimuli/img_2pnl2.png"
$ws.Range("M10").Value = 6.621621621621622
$ws.Range("N10").Value = 7.135135135135135
$ws.Range("O10").Value = 6.878378378378379
$ws.Range("P10").Value = 37
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = 1
$ws.Range("S10").Value = 1
$ws.Range("T10").Value = 1
$ws.Range("U10").Value = 1
$ws.Range("V10").Value = 1

# Row 11
$ws.Range("I11").ClearContents()
$ws.Range("J11").Value = "new"
$ws.Range("K11").Value = "f"
$ws.Range("L11").Value = "stimuli/img_9z99v.png"
$ws.Range("M11").Value = 81.15625
$ws.Range("N11").Value = 64.78125
$ws.Range("O11").Value = 72.96875
$ws.Range("P11").Value = 32
$ws.Range("Q11").Value = 8
$ws.Range("R11").Value = 8
$ws.Range("S11").Value = 8
$ws.Range("T11").Value = 8
$ws.Range("U11").Value = 8
$ws.Range("V11").Value = 8

# Row 12
$ws.Range("I12").ClearContents()
$ws.Range("J12").Value = "new"
$ws.Range("K12").Value = "f"
$ws.Range("L12").Value = "stimuli/img_qbdgm.png"
$ws.Range("M12").Value = 76.88095238095238
$ws.Range("N12").Value = 60.40476190476191
$ws.Range("O12").Value = 68.64285714285714
$ws.Range("P12").Value = 42
$ws.Range("Q12").Value = 7
$ws.Range("R12").Value = 7
$ws.Range("S12").Value = 7
$ws.Range("T12").Value = 7
$ws.Range("U12").Value = 7
$ws.Range("V12").Value = 7

# Row 13
$ws.Range("I13").Value = "target"
$ws.Range("J13").Value = "old"
$ws.Range("K13").Value = "j"
$ws.Range("L13").Value = "stimuli/img_v8dra.png"
$ws.Range("M13").Value = 61.77272727272727
$ws.Range("N13").Value = 38.79545454545455
$ws.Range("O13").Value = 50.28409090909091
$ws.Range("P13").Value = 44
$ws.Range("Q13").Value = 3
$ws.Range("R13").Value = 3
$ws.Range("S13").Value = 3
$ws.Range("T13").Value = 3
$ws.Range("U13").Value = 4
$ws.Range("V13").Value = 3

# Row 14
$ws.Range("I14").ClearContents()
$ws.Range("J14").Value = "new"
$ws.Range("K14").Value = "f"
$ws.Range("L14").Value = "stimuli/img_sltwe.png"
$ws.Range("M14").Value = 72.025
$ws.Range("N14").Value = 46.875
$ws.Range("O14").Value = 59.45
$ws.Range("P14").Value = 40
$ws.Range("Q14").Value = 5
$ws.Range("R14").Value = 5
$ws.Range("S14").Value = 5
$ws.Range("T14").Value = 5
$ws.Range("U14").Value = 5
$ws.Range("V14").Value = 5

# Row 15
$ws.Range("I15").ClearContents()
$ws.Range("J15").Value = "new"
$ws.Range("K15").Value = "f"
$ws.Range("L15").Value = "stimuli/img_dmjh8.png"
$ws.Range("M15").Value = 57.48648648648648
$ws.Range("N15").Value = 37.64864864864865
$ws.Range("O15").Value = 47.56756756756756
$ws.Range("P15").Value = 37
$ws.Range("Q15").Value = 3
$ws.Range("R15").Value = 3
$ws.Range("S15").Value = 3
$ws.Range("T15").Value = 3
$ws.Range("U15").Value = 3
$ws.Range("V15").Value = 3

# Row 16
$ws.Range("I16").ClearContents()
$ws.Range("J16").Value = "new"
$ws.Range("K16").Value = "f"
$ws.Range("L16").Value = "stimuli/img_kn0we.png"
$ws.Range("M16").Value = 80.1590909090909
$ws.Range("N16").Value = 56.68181818181818
$ws.Range("O16").Value = 68.42045454545455
$ws.Range("P16").Value = 44
$ws.Range("Q16").Value = 7
$ws.Range("R16").Value = 7
$ws.Range("S16").Value = 7
$ws.Range("T16").Value = 7
$ws.Range("U16").Value = 7
$ws.Range("V16").Value = 7

# Row 17
$ws.Range("I17").ClearContents()
$ws.Range("J17").Value = "new"
$ws.Range("K17").Value = "f"
$ws.Range("L17").Value = "stimuli/img_7os7q.png"
$ws.Range("M17").Value = 59.7027027027027
$ws.Range("N17").Value = 34.94594594594594
$ws.Range("O17").Value = 47.32432432432432
$ws.Range("P17").Value = 37
$ws.Range("Q17").Value = 3
$ws.Range("R17").Value = 3
$ws.Range("S17").Value = 3
$ws.Range("T17").Value = 3
$ws.Range("U17").Value = 3
$ws.Range("V17").Value = 3

# Row 18
$ws.Range("I18").ClearContents()
$ws.Range("J18").Value = "new"
$ws.Range("K18").Value = "f"
$ws.Range("L18").Value = "stimuli/img_kzg3h.png"
$ws.Range("M18").Value = 77.02777777777777
$ws.Range("N18").Value = 56.22222222222222
$ws.Range("O18").Value = 66.625
$ws.Range("P18").Value = 36
$ws.Range("Q18").Value = 7
$ws.Range("R18").Value = 7
$ws.Range("S18").Value = 7
$ws.Range("T18").Value = 7
$ws.Range("U18").Value = 7
$ws.Range("V18").Value = 7

# Row 19
$ws.Range("I19").ClearContents()
$ws.Range("J19").Value = "new"
$ws.Range("K19").Value = "f"
$ws.Range("L19").Value = "stimuli/img_jge7p.png"
$ws.Range("M19").Value = 90.42424242424242
$ws.Range("N19").Value = 75.63636363636364
$ws.Range("O19").Value = 83.03030303030303
$ws.Range("P19").Value = 33
$ws.Range("Q19").Value = 10
$ws.Range("R19").Value = 10
$ws.Range("S19").Value = 10
$ws.Range("T19").Value = 10
$ws.Range("U19").Value = 10
$ws.Range("V19").Value = 10

# Row 20
$ws.Range("I20").ClearContents()
$ws.Range("J20").Value = "new"
$ws.Range("K20").Value = "f"
$ws.Range("L20").Value = "stimuli/img_o37la.png"
$ws.Range("M20").Value = 65.24324324324324
$ws.Range("N20").Value = 42.78378378378378
$ws.Range("O20").Value = 54.01351351351352
$ws.Range("P20").Value = 37
$ws.Range("Q20").Value = 4
$ws.Range("R20").Value = 4
$ws.Range("S20").Value = 4
$ws.Range("T20").Value = 4
$ws.Range("U20").Value = 4
$ws.Range("V20").Value = 4

# Row 21
$ws.Range("I21").ClearContents()
$ws.Range("J21").Value = "new"
$ws.Range("K21").Value = "f"
$ws.Range("L21").Value = "stimuli/img_65cdi.png"
$ws.Range("M21").Value = 46.92307692307692
$ws.Range("N21").Value = 27
$ws.Range("O21").Value = 36.96153846153846
$ws.Range("P21").Value = 39
$ws.Range("Q21").Value = 2
$ws.Range("R21").Value = 2
$ws.Range("S21").Value = 2
$ws.Range("T21").Value = 2
$ws.Range("U21").Value = 2
$ws.Range("V21").Value = 2

# Row 22
$ws.Range("I22").ClearContents()
$ws.Range("J22").Value = "new"
$ws.Range("K22").Value = "f"
$ws.Range("L22").Value = "stimuli/img_fbihy.png"
$ws.Range("M22").Value = 44.39024390243902
$ws.Range("N22").Value = 26.90243902439024
$ws.Range("O22").Value = 35.64634146341464
$ws.Range("P22").Value = 41
$ws.Range("Q22").Value = 2
$ws.Range("R22").Value = 2
$ws.Range("S22").Value = 2
$ws.Range("T22").Value = 2
$ws.Range("U22").Value = 2
$ws.Range("V22").Value = 2

# Row 23
$ws.Range("I23").Value = "target"
$ws.Range("J23").Value = "old"
$ws.Range("K23").Value = "j"
$ws.Range("L23").Value = "stimuli/img_cgdyc.png"
$ws.Range("M23").Value = 32.93023255813954
$ws.Range("N23").Value = 14.04651162790698
$ws.Range("O23").Value = 23.48837209302326
$ws.Range("P23").Value = 43
$ws.Range("Q23").Value = 1
$ws.Range("R23").Value = 1
$ws.Range("S23").Value = 1
$ws.Range("T23").Value = 1
$ws.Range("U23").Value = 1
$ws.Range("V23").Value = 1

# Row 24
$ws.Range("I24").ClearContents()
$ws.Range("J24").Value = "new"
$ws.Range("K24").Value = "f"
$ws.Range("L24").Value = "stimuli/img_yteqw.png"
$ws.Range("M24").Value = 66.83783783783784
$ws.Range("N24").Value = 43.78378378378378
$ws.Range("O24").Value = 55.31081081081081
$ws.Range("P24").Value = 37
$ws.Range("Q24").Value = 4
$ws.Range("R24").Value = 4
$ws.Range("S24").Value = 4
$ws.Range("T24").Value = 5
$ws.Range("U24").Value = 4
$ws.Range("V24").Value = 4

# Row 25
$ws.Range("I25").Value = "target"
$ws.Range("J25").Value = "old"
$ws.Range("K25").Value = "j"
$ws.Range("L25").Value = "stimuli/img_fqgem.png"
$ws.Range("M25").Value = 80.75
$ws.Range("N25").Value = 61.475
$ws.Range("O25").Value = 71.1125
$ws.Range("P25").Value = 40
$ws.Range("Q25").Value = 8
$ws.Range("R25").Value = 8
$ws.Range("S25").Value = 8
$ws.Range("T25").Value = 8
$ws.Range("U25").Value = 8
$ws.Range("V25").Value = 8

# Row 26
$ws.Range("I26").ClearContents()
$ws.Range("J26").Value = "catch"
$ws.Range("K26").Value = "f"
$ws.Range("L26").Value = "stimuli/catch_05.jpg"
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("O26").ClearContents()
$ws.Range("P26").ClearContents()
$ws.Range("Q26").ClearContents()
$ws.Range("R26").ClearContents()
$ws.Range("S26").ClearContents()
$ws.Range("T26").ClearContents()
$ws.Range("U26").ClearContents()
$ws.Range("V26").ClearContents()

# Row 27
$ws.Range("I27").ClearContents()
$ws.Range("J27").Value = "new"
$ws.Range("K27").Value = "f"
$ws.Range("L27").Value = "stimuli/img_l1h36.png"
$ws.Range("M27").Value = 26.64285714285714
$ws.Range("N27").Value = 9.142857142857142
$ws.Range("O27").Value = 17.89285714285714
$ws.Range("P27").Value = 42
$ws.Range("Q27").Value = 1
$ws.Range("R27").Value = 1
$ws.Range("S27").Value = 1
$ws.Range("T27").Value = 1
$ws.Range("U27").Value = 1
$ws.Range("V27").Value = 1

# Row 28
$ws.Range("I28").ClearContents()
$ws.Range("J28").Value = "new"
$ws.Range("K28").Value = "f"
$ws.Range("L28").Value = "stimuli/img_gbypq.png"
$ws.Range("M28").Value = 76.275
$ws.Range("N28").Value = 51.925
$ws.Range("O28").Value = 64.1
$ws.Range("P28").Value = 40
$ws.Range("Q28").Value = 6
$ws.Range("R28").Value = 6
$ws.Range("S28").Value = 6
$ws.Range("T28").Value = 6
$ws.Range("U28").Value = 6
$ws.Range("V28").Value = 6

# Row 29
$ws.Range("I29").ClearContents()
$ws.Range("J29").Value = "new"
$ws.Range("K29").Value = "f"
$ws.Range("L29").Value = "stimuli/img_zv0dq.png"
$ws.Range("M29").Value = 76.86842105263158
$ws.Range("N29").Value = 52.71052631578947
$ws.Range("O29").Value = 64.78947368421052
$ws.Range("P29").Value = 38
$ws.Range("Q29").Value = 6
$ws.Range("R29").Value = 6
$ws.Range("S29").Value = 6
$ws.Range("T29").Value = 6
$ws.Range("U29").Value = 6
$ws.Range("V29").Value = 6

# Row 30
$ws.Range("I30").ClearContents()
$ws.Range("J30").Value = "new"
$ws.Range("K30").Value = "f"
$ws.Range("L30").Value = "stimuli/img_9pfbj.png"
$ws.Range("M30").Value = 91.27272727272727
$ws.Range("N30").Value = 80.0909090909091
$ws.Range("O30").Value = 85.68181818181819
$ws.Range("P30").Value = 33
$ws.Range("Q30").Value = 10
$ws.Range("R30").Value = 10
$ws.Range("S30").Value = 10
$ws.Range("T30").Value = 10
$ws.Range("U30").Value = 10
$ws.Range("V30").Value = 10

# Row 31
$ws.Range("I31").ClearContents()
$ws.Range("J31").Value = "new"
$ws.Range("K31").Value = "f"
$ws.Range("L31").Value = "stimuli/img_sfh4b.png"
$ws.Range("M31").Value = 69.06521739130434
$ws.Range("N31").Value = 49.54347826086956
$ws.Range("O31").Value = 59.30434782608695
$ws.Range("P31").Value = 46
$ws.Range("Q31").Value = 5
$ws.Range("R31").Value = 5
$ws.Range("S31").Value = 5
$ws.Range("T31").Value = 5
$ws.Range("U31").Value = 5
$ws.Range("V31").Value = 5

# Row 32
$ws.Range("I32").Value = "target"
$ws.Range("J32").Value = "old"
$ws.Range("K32").Value = "j"
$ws.Range("L32").Value = "stimuli/img_i7vab.png"
$ws.Range("M32").Value = 86.4
$ws.Range("N32").Value = 67.8
$ws.Range("O32").Value = 77.1
$ws.Range("P32").Value = 35
$ws.Range("Q32").Value = 9
$ws.Range("R32").Value = 9
$ws.Range("S32").Value = 9
$ws.Range("T32").Value = 9
$ws.Range("U32").Value = 9
$ws.Range("V32").Value = 9

# Row 33
$ws.Range("I33").ClearContents()
$ws.Range("J33").Value = "new"
$ws.Range("K33").Value = "f"
$ws.Range("L33").Value = "stimuli/img_badai.png"
$ws.Range("M33").Value = 63.97435897435897
$ws.Range("N33").Value = 43.38461538461539
$ws.Range("O33").Value = 53.67948717948718
$ws.Range("P33").Value = 39
$ws.Range("Q33").Value = 4
$ws.Range("R33").Value = 4
$ws.Range("S33").Value = 4
$ws.Range("T33").Value = 4
$ws.Range("U33").Value = 4
$ws.Range("V33").Value = 4

# Row 34
$ws.Range("I34").Value = "target"
$ws.Range("J34").Value = "old"
$ws.Range("K34").Value = "j"
$ws.Range("L34").Value = "stimuli/img_aweye.png"
$ws.Range("M34").Value = 53.42105263157895
$ws.Range("N34").Value = 31.84210526315789
$ws.Range("O34").Value = 42.63157894736842
$ws.Range("P34").Value = 38
$ws.Range("Q34").Value = 2
$ws.Range("R34").Value = 2
$ws.Range("S34").Value = 2
$ws.Range("T34").Value = 3
$ws.Range("U34").Value = 3
$ws.Range("V34").Value = 2

# Row 35
$ws.Range("I35").ClearContents()
$ws.Range("J35").Value = "new"
$ws.Range("K35").Value = "f"
$ws.Range("L35").Value = "stimuli/img_c2pbs.png"
$ws.Range("M35").Value = 21.95238095238095
$ws.Range("N35").Value = 14.47619047619048
$ws.Range("O35").Value = 18.21428571428572
$ws.Range("P35").Value = 42
$ws.Range("Q35").Value = 1
$ws.Range("R35").Value = 1
$ws.Range("S35").Value = 1
$ws.Range("T35").Value = 1
$ws.Range("U35").Value = 1
$ws.Range("V35").Value = 1

# Row 36
$ws.Range("I36").ClearContents()
$ws.Range("J36").Value = "new"
$ws.Range("K36").Value = "f"
$ws.Range("L36").Value = "stimuli/img_bntrh.png"
$ws.Range("M36").Value = 76.07894736842105
$ws.Range("N36").Value = 53.36842105263158
$ws.Range("O36").Value = 64.72368421052632
$ws.Range("P36").Value = 38
$ws.Range("Q36").Value = 6
$ws.Range("R36").Value = 6
$ws.Range("S36").Value = 6
$ws.Range("T36").Value = 6
$ws.Range("U36").Value = 6
$ws.Range("V36").Value = 6

# Row 37
$ws.Range("I37").Value = "target"
$ws.Range("J37").Value = "old"
$ws.Range("K37").Value = "j"
$ws.Range("L37").Value = "stimuli/img_3bxjb.png"
$ws.Range("M37").Value = 87.28571428571429
$ws.Range("N37").Value = 72.65714285714286
$ws.Range("O37").Value = 79.97142857142858
$ws.Range("P37").Value = 35
$ws.Range("Q37").Value = 10
$ws.Range("R37").Value = 10
$ws.Range("S37").Value = 10
$ws.Range("T37").Value = 9
$ws.Range("U37").Value = 9
$ws.Range("V37").Value = 10

# Row 38
$ws.Range("I38").ClearContents()
$ws.Range("J38").Value = "new"
$ws.Range("K38").Value = "f"
$ws.Range("L38").Value = "stimuli/img_anzgh.png"
$ws.Range("M38").Value = 75.10526315789474
$ws.Range("N38").Value = 55.76315789473684
$ws.Range("O38").Value = 65.4342105263158
$ws.Range("P38").Value = 38
$ws.Range("Q38").Value = 6
$ws.Range("R38").Value = 6
$ws.Range("S38").Value = 6
$ws.Range("T38").Value = 6
$ws.Range("U38").Value = 6
$ws.Range("V38").Value = 6

# Row 39
$ws.Range("I39").Value = "target"
$ws.Range("J39").Value = "old"
$ws.Range("K39").Value = "j"
$ws.Range("L39").Value = "stimuli/img_ose78.png"
$ws.Range("M39").Value = 80.19444444444444
$ws.Range("N39").Value = 60.25
$ws.Range("O39").Value = 70.22222222222223
$ws.Range("P39").Value = 36
$ws.Range("Q39").Value = 8
$ws.Range("R39").Value = 7
$ws.Range("S39").Value = 7
$ws.Range("T39").Value = 7
$ws.Range("U39").Value = 7
$ws.Range("V39").Value = 7

# Row 40
$ws.Range("I40").Value = "target"
$ws.Range("J40").Value = "old"
$ws.Range("K40").Value = "j"
$ws.Range("L40").Value = "stimuli/img_ic3os.png"
$ws.Range("M40").Value = 84.79069767441861
$ws.Range("N40").Value = 66.16279069767442
$ws.Range("O40").Value = 75.47674418604652
$ws.Range("P40").Value = 43
$ws.Range("Q40").Value = 9
$ws.Range("R40").Value = 9
$ws.Range("S40").Value = 9
$ws.Range("T40").Value = 8
$ws.Range("U40").Value = 9
$ws.Range("V40").Value = 9

# Row 41
$ws.Range("I41").ClearContents()
$ws.Range("J41").Value = "new"
$ws.Range("K41").Value = "f"
$ws.Range("L41").Value = "stimuli/img_5p2ql.png"
$ws.Range("M41").Value = 89.19565217391305
$ws.Range("N41").Value = 72.52173913043478
$ws.Range("O41").Value = 80.8586956521739
$ws.Range("P41").Value = 46
$ws.Range("Q41").Value = 10
$ws.Range("R41").Value = 10
$ws.Range("S41").Value = 10
$ws.Range("T41").Value = 10
$ws.Range("U41").Value = 10
$ws.Range("V41").Value = 9

# Row 42
$ws.Range("I42").ClearContents()
$ws.Range("J42").Value = "new"
$ws.Range("K42").Value = "f"
$ws.Range("L42").Value = "stimuli/img_t4hvr.png"
$ws.Range("M42").Value = 61.69230769230769
$ws.Range("N42").Value = 39.76923076923077
$ws.Range("O42").Value = 50.73076923076923
$ws.Range("P42").Value = 39
$ws.Range("Q42").Value = 3
$ws.Range("R42").Value = 3
$ws.Range("S42").Value = 3
$ws.Range("T42").Value = 4
$ws.Range("U42").Value = 3
$ws.Range("V42").Value = 4
